# Generate Report for Handback
# The handback status report is regenerated: the "in sync with en-US" status
# becomes "not in sync with en-US", and the zh-cn / de-de handback datetimes
# for the first file (1169e9a6-...) are refreshed to new timestamps.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: not in sync with en-US"

# --- Overview sheet: zh-cn / de-de status columns (E, F) for both rows ---
$wsOverview.Range("E2").Value2 = $newStatus
$wsOverview.Range("F2").Value2 = $newStatus
$wsOverview.Range("E3").Value2 = $newStatus
$wsOverview.Range("F3").Value2 = $newStatus

# --- zh-cn detail sheet ---
$wsZhCn.Range("C2").Value2 = $newStatus
$wsZhCn.Range("C3").Value2 = $newStatus
# Correspond Handback DateTime refreshed for the first file's handback
$wsZhCn.Range("K2").Value2 = "2016-09-06 15:51:47"

# --- de-de detail sheet ---
$wsDeDe.Range("C2").Value2 = $newStatus
$wsDeDe.Range("C3").Value2 = $newStatus
# Correspond Handback DateTime refreshed for the first file's handback
$wsDeDe.Range("K2").Value2 = "2016-09-06 15:52:15"

# --- Column widths widened to fit the longer status text ---
$wsOverview.Columns.Item(5).ColumnWidth = 32.6
$wsOverview.Columns.Item(6).ColumnWidth = 32.6
$wsZhCn.Columns.Item(3).ColumnWidth = 32.6
$wsDeDe.Columns.Item(3).ColumnWidth = 32.6
